# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used data row (header is row 1, data starts row 2).
$lastRow = $ws.UsedRange.Rows.Count

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting from an existing header cell (AC1) onto the new headers
# so they match the bold/bordered/centered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows: AD = 86 (Wins), AE = 76 (Losses), AF = 0 (Ties) ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}

$ws.Range("A1").Select()
